# Apply the table style change on the financial-documents table (slide 5).
#
# The table currently uses the deck's custom "Table_0" style
# ({3CACF4BB-14DE-423B-943C-D5984CF458D3}); the edit switches it to the
# built-in table style {26527EDE-59EB-4736-8D4E-5B47A16CBB89}.
#
# Table styles cannot be reassigned through the Table.Style property
# (PowerPoint COM raises an error telling you to use ApplyStyle instead),
# so we call Table.ApplyStyle with the target style id.

$p = $ppt.ActivePresentation

$targetStyleId = "{26527EDE-59EB-4736-8D4E-5B47A16CBB89}"
$applied = $false

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle($targetStyleId)
            $applied = $true
        }
    }
}

if (-not $applied) {
    Write-Host "Warning: no table shape found to restyle"
}
